$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 433; this shifts the existing rows
# 433-482 down to 434-483 (matching the dimension growing from
# A1:R482 to A1:R483).
$ws.Rows(433).Insert()

# Populate the newly inserted row 433 with the new weekly price record.
$ws.Range("A433").Value2 = 9
$ws.Range("B433").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C433").Value2 = "Metropolitana"
$ws.Range("D433").Value2 = 44918
$ws.Range("E433").Value2 = 13
$ws.Range("F433").Value2 = 100112044
$ws.Range("G433").Value2 = "Perejil"
$ws.Range("H433").Value2 = "Sin especificar"
$ws.Range("I433").Value2 = "Primera"
$ws.Range("J433").Value2 = 70
$ws.Range("K433").Value2 = 12000
$ws.Range("L433").Value2 = 12000
$ws.Range("M433").Value2 = 12000
$ws.Range("N433").Value2 = "$/docena de atados"
$ws.Range("O433").Value2 = "Región Metropolitana"
$ws.Range("P433").Value2 = 4000
$ws.Range("Q433").Value2 = 3
$ws.Range("R433").Value2 = "Hortaliza"
